{"js": "// The document contains a \"use case\" summary table whose first cell\n// holds the use-case code (\"CU21\"). The edit renames that code to\n// \"CU20\" (the rest of the table \u2014 description, actors, etc. \u2014 is\n// unchanged).\nconst results = context.document.body.search(\"CU21\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"CU20\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The first table's header cell holds the use-case code \"CU21\".\n# Rename it to \"CU20\" (the rest of the row/table is untouched).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"CU21\"\n$find.Replacement.Text = \"CU20\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
